# This script reorders the per-trial columns (category, condition, stimulus,
# correct answer and the associated rating/statistics columns) across the
# data rows of the sheet, according to a fixed row permutation. Columns
# A-G and J (subject/task/block bookkeeping columns, which are identical
# across all rows or are simple running counters) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (1-based worksheet row numbers).
# The value that ends up in destination row R (for the columns listed below)
# is the value that used to be in source row perm[R] before the edit.
$perm = @{
    2  = 37; 3  = 34; 4  = 24; 5  = 14; 6  = 38; 7  = 5;  8  = 12; 9  = 25;
    10 = 29; 11 = 19; 12 = 4;  13 = 35; 14 = 22; 15 = 16; 16 = 39; 17 = 33;
    18 = 3;  19 = 17; 20 = 10; 21 = 6;  22 = 31; 23 = 28; 24 = 9;  25 = 30;
    26 = 40; 27 = 2;  28 = 26; 29 = 32; 30 = 21; 31 = 7;  32 = 15; 33 = 20;
    34 = 8;  35 = 11; 36 = 23; 37 = 27; 38 = 18; 39 = 41; 40 = 36; 41 = 13
}

# Columns that move together as a unit when rows are permuted.
$cols = @(8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)
# 8=H(category) 9=I(cond_cat) 11=K(correct_answer) 12=L(stimulus)
# 13=M(conceptual) 14=N(perceptual) 15=O(typicality) 16=P(n)
# 17=Q(p_typicality) 18=R(p_conceptual) 19=S(p_perceptual)
# 20=T(r_typicality) 21=U(r_conceptual) 22=V(r_perceptual)

$firstRow = 2
$lastRow = 41

# Snapshot the original values for the affected columns before writing
# anything back, since several source/destination rows overlap (the
# permutation contains multi-row cycles).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write back the permuted values.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $perm[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
}
